# Applies the cryptos-list price/volume update described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.222.12'
$ws.Range('E2').Value = '  -4.75%  '
$ws.Range('D3').Value = '3.085.70'
$ws.Range('E3').Value = '  -5.29%  '
$ws.Range('E4').Value = '  +0.13%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '548.12'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  -5.07%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.80'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  -11.40%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = '3.072.84'
$ws.Range('E8').Value = '  -5.36%  '
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.496'
$ws.Range('D9').Style = $origStyle
$ws.Range('E9').Value = '  -3.82%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.155'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  -6.08%  '
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.22'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  -11.73%  '
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.468'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  -4.65%  '
$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '35.41'
$ws.Range('D13').Style = $origStyle
$ws.Range('E13').Value = '  -6.62%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000216'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  -8.82%  '
$ws.Range('D15').Value = '3.584.97'
$ws.Range('E15').Value = '  -5.19%  '
$ws.Range('D16').Value = '63.262.95'
$ws.Range('E16').Value = '  -4.74%  '
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.111'
$ws.Range('D17').Style = $origStyle
$ws.Range('E17').Value = '  -3.41%  '
$ws.Range('D18').Value = '3.092.98'
$ws.Range('E18').Value = '  -5.13%  '
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.71'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  -5.98%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '487.89'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  -13.12%  '
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.63'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  -5.93%  '
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.717'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  -3.79%  '
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.25'
$ws.Range('D23').Style = $origStyle
$ws.Range('E23').Value = '  -7.43%  '
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.86'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  -4.24%  '
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.36'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  -9.52%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.76'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  -7.61%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.44'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  -10.08%  '
$ws.Range('B29').Value = 'FirstDigitalUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.97'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  -12.16%  '
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.58'
$ws.Range('D31').Style = $origStyle
$ws.Range('E31').Value = '  -4.95%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.13'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  -4.37%  '
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.50'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  -9.85%  '
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '58.05'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  +5.10%  '
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '513.93'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  -9.12%  '
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.02'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  -6.23%  '
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.11'
$ws.Range('D37').Style = $origStyle
$ws.Range('E37').Value = '  -11.55%  '
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0401'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  -12.99%  '
$ws.Range('D39').Value = '3.144.59'
$ws.Range('E39').Value = '  -0.60%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0802'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  -7.91%  '
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.119'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  -7.52%  '
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.15'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  -5.77%  '
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.65'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  -15.02%  '
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.259'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  -6.16%  '
$ws.Range('E46').Value = '  -10.66%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '25.14'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  -5.68%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '120.97'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  -3.27%  '
$ws.Range('D50').Value = '0.0₃0502'
$ws.Range('E50').Value = '  -10.45%  '
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.33'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  +30.16%  '
